$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) stays as text so values like trailing zeros are preserved
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "28.083.25"
$ws.Range("E2").Value = "  +0.17%  "

$ws.Range("D3").Value = "1.817.08"
$ws.Range("E3").Value = "  +2.26%  "

$ws.Range("E4").Value = "  -0.33%  "

$ws.Range("D5").Value = "338.19"
$ws.Range("E5").Value = "  -0.35%  "

$ws.Range("D6").Value = "0.9995"
$ws.Range("E6").Value = "  -0.13%  "

$ws.Range("D7").Value = "0.4290"
$ws.Range("E7").Value = "  +12.40%  "

$ws.Range("D8").Value = "0.3508"
$ws.Range("E8").Value = "  +2.61%  "

$ws.Range("D9").Value = "45.78"
$ws.Range("E9").Value = "  -2.34%  "

$ws.Range("D10").Value = "1.150"
$ws.Range("E10").Value = "  +0.77%  "

$ws.Range("D11").Value = "0.07447"
$ws.Range("E11").Value = "  +0.77%  "

$ws.Range("D12").Value = "22.99"
$ws.Range("E12").Value = "  -1.30%  "

$ws.Range("D13").Value = "1.001"
$ws.Range("E13").Value = "  -0.20%  "

$ws.Range("D14").Value = "6.260"
$ws.Range("E14").Value = "  -1.90%  "

$ws.Range("D15").Value = "1.815.12"
$ws.Range("E15").Value = "  +2.11%  "

$ws.Range("D16").Value = "7.275"
$ws.Range("E16").Value = "  -1.87%  "

$ws.Range("E17").Value = "  +1.01%  "

$ws.Range("D18").Value = "0.06675"
$ws.Range("E18").Value = "  +0.24%  "

$ws.Range("D19").Value = "81.99"
$ws.Range("E19").Value = "  -0.31%  "

$ws.Range("D20").Value = "0.9998"
$ws.Range("E20").Value = "  -0.14%  "

$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "6.470"
$ws.Range("E21").Value = "  +0.94%  "

$ws.Range("B22").Value = "Avalanche"
$ws.Range("C22").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D22").Value = "17.28"
$ws.Range("E22").Value = "  -0.45%  "

$ws.Range("D23").Value = "28.097.72"
$ws.Range("E23").Value = "  +0.06%  "

$ws.Range("D24").Value = "12.02"
$ws.Range("E24").Value = "  -0.44%  "

$ws.Range("D25").Value = "2.383"
$ws.Range("E25").Value = "  +0.13%  "

$ws.Range("D26").Value = "2.496"

$ws.Range("D27").Value = "20.72"
$ws.Range("E27").Value = "  +0.09%  "

$ws.Range("D28").Value = "156.14"
$ws.Range("E28").Value = "  +1.19%  "

$ws.Range("D29").Value = "2.023.26"
$ws.Range("E29").Value = "  +2.16%  "

$ws.Range("D30").Value = "1.303"
$ws.Range("E30").Value = "  -10.07%  "

$ws.Range("D31").Value = "132.60"
$ws.Range("E31").Value = "  -1.12%  "

$ws.Range("D32").Value = "4.058"
$ws.Range("E32").Value = "  +0.84%  "

$ws.Range("D33").Value = "5.964"
$ws.Range("E33").Value = "  -1.60%  "

$ws.Range("D34").Value = "0.09213"
$ws.Range("E34").Value = "  +3.73%  "

$ws.Range("D35").Value = "12.37"
$ws.Range("E35").Value = "  -2.71%  "

$ws.Range("D36").Value = "0.02365"
$ws.Range("E36").Value = "  -1.30%  "

$ws.Range("D37").Value = "0.6738"
$ws.Range("E37").Value = "  -1.46%  "

$ws.Range("D38").Value = "5.245"
$ws.Range("E38").Value = "  -0.65%  "

$ws.Range("D39").Value = "0.06265"
$ws.Range("E39").Value = "  -2.25%  "

$ws.Range("D40").Value = "0.2166"
$ws.Range("E40").Value = "  +0.10%  "

$ws.Range("D41").Value = "1.495"
$ws.Range("E41").Value = "  -0.23%  "

$ws.Range("D42").Value = "1.217"
$ws.Range("E42").Value = "  -1.54%  "

$ws.Range("D43").Value = "8.228"
$ws.Range("E43").Value = "  +0.12%  "

$ws.Range("D44").Value = "0.9985"
$ws.Range("E44").Value = "  -0.24%  "

$ws.Range("D45").Value = "14.07"
$ws.Range("E45").Value = "  -1.36%  "

$ws.Range("D46").Value = "3.871"
$ws.Range("E46").Value = "  +0.15%  "

$ws.Range("D47").Value = "0.6130"
$ws.Range("E47").Value = "  -2.14%  "

$ws.Range("D48").Value = "128.71"
$ws.Range("E48").Value = "  -3.09%  "

$ws.Range("E49").Value = "  -0.98%  "

$ws.Range("D50").Value = "1.180"
$ws.Range("E50").Value = "  -3.09%  "

$ws.Range("D51").Value = "0.07109"
$ws.Range("E51").Value = "  -5.20%  "
